# Update the "想去人数" (want-to-go count) values in column F
# for both the "展览" and "全部类型" worksheets, which carry the
# same data table.

$wb = $excel.ActiveWorkbook

# Row number -> new value for column F
$updates = @{
    3  = 3410
    4  = 250
    5  = 148
    6  = 222
    7  = 1747
    8  = 1660
    9  = 477
    10 = 382
    12 = 33
    14 = 41
    15 = 236
    16 = 7
    17 = 40
    23 = 64
    24 = 57
    26 = 412
    27 = 276
    32 = 485
    33 = 2333
    36 = 486
    37 = 580
    43 = 549
    44 = 417
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
